$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "'69.951.64"
Set-TextValue $ws.Range("E2") "'  -1.13%  "

# Row 3
Set-TextValue $ws.Range("D3") "'2.511.02"
Set-TextValue $ws.Range("E3") "'  -2.30%  "

# Row 4
Set-TextValue $ws.Range("E4") "'  +0.10%  "

# Row 5
Set-TextValue $ws.Range("D5") "'572.81"
Set-TextValue $ws.Range("E5") "'  -0.83%  "

# Row 6
Set-TextValue $ws.Range("D6") "'166.80"
Set-TextValue $ws.Range("E6") "'  -2.32%  "

# Row 7
Set-TextValue $ws.Range("E7") "'  +0.01%  "

# Row 8
Set-TextValue $ws.Range("D8") "'0.519"
Set-TextValue $ws.Range("E8") "'  +1.39%  "

# Row 9
Set-TextValue $ws.Range("D9") "'2.512.28"
Set-TextValue $ws.Range("E9") "'  -2.21%  "

# Row 10
Set-TextValue $ws.Range("D10") "'0.161"
Set-TextValue $ws.Range("E10") "'  -3.45%  "

# Row 11
Set-TextValue $ws.Range("D11") "'0.167"
Set-TextValue $ws.Range("E11") "'  -0.88%  "

# Row 12
Set-TextValue $ws.Range("D12") "'0.355"
Set-TextValue $ws.Range("E12") "'  +1.62%  "

# Row 13
Set-TextValue $ws.Range("D13") "'4.90"
Set-TextValue $ws.Range("E13") "'  +0.99%  "

# Row 14
Set-TextValue $ws.Range("D14") "'2.982.72"
Set-TextValue $ws.Range("E14") "'  -1.94%  "

# Row 15
Set-TextValue $ws.Range("D15") "'70.020.47"
Set-TextValue $ws.Range("E15") "'  -0.92%  "

# Row 16
Set-TextValue $ws.Range("D16") "'0.0000177"
Set-TextValue $ws.Range("E16") "'  -3.45%  "

# Row 17
Set-TextValue $ws.Range("D17") "'25.00"
Set-TextValue $ws.Range("E17") "'  -1.08%  "

# Row 18
Set-TextValue $ws.Range("D18") "'2.525.04"
Set-TextValue $ws.Range("E18") "'  -1.78%  "

# Row 19
Set-TextValue $ws.Range("D19") "'7.86"
Set-TextValue $ws.Range("E19") "'  +6.25%  "

# Row 20
Set-TextValue $ws.Range("D20") "'11.34"
Set-TextValue $ws.Range("E20") "'  -4.02%  "

# Row 21
Set-TextValue $ws.Range("D21") "'349.26"
Set-TextValue $ws.Range("E21") "'  -4.04%  "

# Row 22
Set-TextValue $ws.Range("D22") "'3.91"
Set-TextValue $ws.Range("E22") "'  -2.22%  "

# Row 23
Set-TextValue $ws.Range("D23") "'1.99"
Set-TextValue $ws.Range("E23") "'  -1.43%  "

# Row 24
Set-TextValue $ws.Range("E24") "'  -0.06%  "

# Row 25
Set-TextValue $ws.Range("D25") "'70.11"
Set-TextValue $ws.Range("E25") "'  -0.21%  "

# Row 26
Set-TextValue $ws.Range("D26") "'3.99"
Set-TextValue $ws.Range("E26") "'  -3.72%  "

# Row 27
Set-TextValue $ws.Range("B27") "'WrappedeETH"
Set-TextValue $ws.Range("C27") "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D27") "'2.694.39"
Set-TextValue $ws.Range("E27") "'  -0.27%  "

# Row 28
Set-TextValue $ws.Range("B28") "'Aptos"
Set-TextValue $ws.Range("C28") "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D28") "'8.85"
Set-TextValue $ws.Range("E28") "'  -5.56%  "

# Row 29
Set-TextValue $ws.Range("D29") "'0.999"
Set-TextValue $ws.Range("E29") "'  +0.02%  "

# Row 30
Set-TextValue $ws.Range("D30") "'0.0₃0904"
Set-TextValue $ws.Range("E30") "'  -3.17%  "

# Row 31
Set-TextValue $ws.Range("D31") "'7.86"
Set-TextValue $ws.Range("E31") "'  +0.29%  "

# Row 32
Set-TextValue $ws.Range("B32") "'Fetch.AI"
Set-TextValue $ws.Range("C32") "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D32") "'1.25"
Set-TextValue $ws.Range("E32") "'  -4.20%  "

# Row 33
Set-TextValue $ws.Range("B33") "'Bittensor"
Set-TextValue $ws.Range("C33") "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D33") "'461.87"
Set-TextValue $ws.Range("E33") "'  -5.17%  "

# Row 34
Set-TextValue $ws.Range("D34") "'1.74"
Set-TextValue $ws.Range("E34") "'  -1.96%  "

# Row 35
Set-TextValue $ws.Range("E35") "'  +0.05%  "

# Row 36
Set-TextValue $ws.Range("D36") "'0.117"
Set-TextValue $ws.Range("E36") "'  +2.58%  "

# Row 37
Set-TextValue $ws.Range("D37") "'156.46"
Set-TextValue $ws.Range("E37") "'  -0.92%  "

# Row 38
Set-TextValue $ws.Range("D38") "'19.06"
Set-TextValue $ws.Range("E38") "'  +1.12%  "

# Row 39
Set-TextValue $ws.Range("D39") "'18.64"
Set-TextValue $ws.Range("E39") "'  -0.85%  "

# Row 40
Set-TextValue $ws.Range("E40") "'  -0.03%  "

# Row 41
Set-TextValue $ws.Range("D41") "'4.77"
Set-TextValue $ws.Range("E41") "'  -0.58%  "

# Row 42
Set-TextValue $ws.Range("D42") "'0.317"
Set-TextValue $ws.Range("E42") "'  -1.55%  "

# Row 43
Set-TextValue $ws.Range("D43") "'1.60"
Set-TextValue $ws.Range("E43") "'  -5.63%  "

# Row 44
Set-TextValue $ws.Range("D44") "'2.31"
Set-TextValue $ws.Range("E44") "'  -7.35%  "

# Row 45
Set-TextValue $ws.Range("B45") "'OKB"
Set-TextValue $ws.Range("C45") "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D45") "'38.33"
Set-TextValue $ws.Range("E45") "'  -0.52%  "

# Row 46
Set-TextValue $ws.Range("B46") "'ImmutableX"
Set-TextValue $ws.Range("C46") "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D46") "'1.15"
Set-TextValue $ws.Range("E46") "'  -14.12%  "

# Row 47
Set-TextValue $ws.Range("D47") "'142.20"
Set-TextValue $ws.Range("E47") "'  -2.85%  "

# Row 48
Set-TextValue $ws.Range("D48") "'0.526"
Set-TextValue $ws.Range("E48") "'  -1.62%  "

# Row 49
Set-TextValue $ws.Range("D49") "'3.48"
Set-TextValue $ws.Range("E49") "'  -2.85%  "

# Row 50
Set-TextValue $ws.Range("D50") "'1.58"
Set-TextValue $ws.Range("E50") "'  -4.10%  "

# Row 51
Set-TextValue $ws.Range("D51") "'0.0730"
Set-TextValue $ws.Range("E51") "'  -0.93%  "

